# Commit: "Modulus opdateret til Modulus Social. Vena har fået status
# godkendt for CPD-DK og XDS Metadata" touches several workbooks in the
# repo; for THIS workbook (Børnejounalsystemer) the only substantive,
# automatable change is the weekly re-dating of the single data sheet -
# the tab is named after the "last updated" date, and the sheet-scoped
# defined name that spans the table must follow it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Opdateret d. 02-12-2025" -> "Opdateret d. 05-12-2025".
# Renaming the sheet automatically repoints the workbook-scoped defined
# name "Børnejounalsystemer" (='Opdateret d. 02-12-2025'!$A$1:$F$62) at
# the new sheet name as well.
$ws.Name = "Opdateret d. 05-12-2025"
